# 16/06/2018 MAMATHA CHICK IN
#
# 1) Merge the two runs that together spell out the WED Jun 13 timestamp
#    into a single run (cosmetic run-merge, same visible text).
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "WED Jun 13 12:17:18 IST 2018", $false, $false, $false, $false, $false,
    $true, 1, $false, "WED Jun 13 12:17:18 IST 2018", 2) | Out-Null

# 2) Append a brand-new purchase entry (FRI Jun 15 / JAYAKKA / CARROT /
#    5922) right after the previous entry's "Amount balance" line, before
#    the block of trailing blank paragraphs that closes the document.
#
# Locate the last "Amount balance" paragraph (the one belonging to the
# "WED Jun 13" entry) so the new block is inserted in the right spot
# regardless of exact paragraph numbering.
$target = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Amount balance*") {
        $target = $i
    }
}

# Anchor = first paragraph right after that "Amount balance" line; every
# new paragraph is inserted immediately before this anchor, so the anchor
# index simply advances by one after each insertion.
$anchorIndex = $target + 1

# -- blank bold paragraph mark --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Font.Bold = 1
$anchorIndex = $anchorIndex + 1

# -- timestamp line --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "FRI Jun 15 11:28:50 IST 2018"
$anchorIndex = $anchorIndex + 1

# -- Person Name --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Person Name`t`t`t`t- JAYAKKA"
$anchorIndex = $anchorIndex + 1

# -- Bill number --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Bill number`t`t`t`t- 5922"
$anchorIndex = $anchorIndex + 1

# -- separator --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "---------------------------------------------------------------"
$anchorIndex = $anchorIndex + 1

# -- Item Name (1st item) --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Item Name`t`t`t`t- CARROT"
$anchorIndex = $anchorIndex + 1

# -- Amount Received (red) --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Font.Color = 255
$rng.Text = "Amount Received`t`t`t- 1720"
$anchorIndex = $anchorIndex + 1

# -- Amount Received mode --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Amount Received mode`t`t- CASH AND CLEARD"
$anchorIndex = $anchorIndex + 1

# -- blank separator paragraph --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$anchorIndex = $anchorIndex + 1

# -- Item Name (2nd item) --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Item Name`t`t`t`t- CARROT"
$anchorIndex = $anchorIndex + 1

# -- Number of Pockets --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Number of Pockets`t`t`t- 2"
$anchorIndex = $anchorIndex + 1

# -- Number of KGs --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Number of KGs`t`t`t- 192"
$anchorIndex = $anchorIndex + 1

# -- Rate --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Rate`t`t`t`t`t- 22"
$anchorIndex = $anchorIndex + 1

# -- Total Price --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Text = "Total Price`t`t`t`t- 4224.0"
$anchorIndex = $anchorIndex + 1

# -- Amount balance (bold) --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Font.Bold = 1
$rng.Text = "Amount balance`t`t`t- 4224.0"
$anchorIndex = $anchorIndex + 1

# -- blank paragraph --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$anchorIndex = $anchorIndex + 1

# -- blank bold paragraph mark (trailing) --
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$rng = $d.Paragraphs($anchorIndex).Range
$rng.Font.Name = "Courier New"
$rng.Font.Bold = 1
$anchorIndex = $anchorIndex + 1

Write-Host "Inserted new entry block starting at paragraph" ($target + 1)
